# Apply the "Tables and Charts correctly postioned" edit:
#  1) Bump the cached datetimeFigureOut footer date from 2025/1/24 to
#     2025/1/25 on the slide master and every slide layout.
#  2) Resize / re-font the "[Vendor Name]" title textbox on slide 1.

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text on the slide master ----------------------
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "2025/1/24") {
            $sh.TextFrame.TextRange.Text = "2025/1/25"
        }
    }
}

# --- 1b) Date placeholder text on every slide layout --------------------
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "2025/1/24") {
                $sh.TextFrame.TextRange.Text = "2025/1/25"
            }
        }
    }
}

# --- 2) Resize / re-font the vendor-name textbox on slide 1 ------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "[Vendor Name]") {
            $sh.Height = 41.198425196850394
            $sh.TextFrame.TextRange.Font.Size = 28
        }
    }
}
